$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price / volume change columns),
# matching the GitHub Actions scheduled data refresh.

$ws.Range('D2').Value = "'37.296.93"
$ws.Range('E2').Value = "'  +1.90%  "
$ws.Range('D3').Value = "'2.083.66"
$ws.Range('E3').Value = "'  -0.40%  "
$ws.Range('D4').Value = "'1.01"
$ws.Range('E4').Value = "'  +0.47%  "
$ws.Range('D5').Value = "'251.40"
$ws.Range('E5').Value = "'  +1.48%  "
$ws.Range('D6').Value = "'0.663"
$ws.Range('E6').Value = "'  -0.42%  "
$ws.Range('E7').Value = "'  +0.14%  "
$ws.Range('D8').Value = "'54.72"
$ws.Range('E8').Value = "'  +21.18%  "
$ws.Range('D9').Value = "'62.26"
$ws.Range('E9').Value = "'  +2.85%  "
$ws.Range('E10').Value = "'  +4.60%  "
$ws.Range('E11').Value = "'  +4.30%  "
$ws.Range('E12').Value = "'  +7.38%  "
$ws.Range('D13').Value = "'15.40"
$ws.Range('E13').Value = "'  +5.64%  "
$ws.Range('B14').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('C14').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D14').Value = "'2.393.07"
$ws.Range('E14').Value = "'  +0.01%  "
$ws.Range('B15').Value = "'Polygon"
$ws.Range('C15').Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range('D15').Value = "'0.860"
$ws.Range('E15').Value = "'  +3.59%  "
$ws.Range('D16').Value = "'5.26"
$ws.Range('E16').Value = "'  +6.55%  "
$ws.Range('D17').Value = "'2.090.57"
$ws.Range('E17').Value = "'  -0.11%  "
$ws.Range('D18').Value = "'37.289.39"
$ws.Range('E18').Value = "'  +1.88%  "
$ws.Range('D19').Value = "'73.20"
$ws.Range('E19').Value = "'  +1.84%  "
$ws.Range('D20').Value = "'14.51"
$ws.Range('E20').Value = "'  +13.53%  "
$ws.Range('D21').Value = "'0.0₃0852"
$ws.Range('E21').Value = "'  +4.55%  "
$ws.Range('D22').Value = "'240.68"
$ws.Range('E22').Value = "'  +0.70%  "
$ws.Range('D23').Value = "'5.26"
$ws.Range('E23').Value = "'  +6.09%  "
$ws.Range('E24').Value = "'  -0.05%  "
$ws.Range('E25').Value = "'  +0.05%  "
$ws.Range('D26').Value = "'171.58"
$ws.Range('E26').Value = "'  +1.32%  "
$ws.Range('D27').Value = "'9.25"
$ws.Range('E27').Value = "'  +4.61%  "
$ws.Range('D28').Value = "'20.83"
$ws.Range('E28').Value = "'  +0.49%  "
$ws.Range('D29').Value = "'2.03"
$ws.Range('E29').Value = "'  +3.00%  "
$ws.Range('E30').Value = "'  +1.78%  "
$ws.Range('D31').Value = "'23.53"
$ws.Range('E31').Value = "'  +6.24%  "
$ws.Range('D32').Value = "'1.09"
$ws.Range('E32').Value = "'  +21.16%  "
$ws.Range('D33').Value = "'4.54"
$ws.Range('E33').Value = "'  +3.09%  "
$ws.Range('D34').Value = "'0.0629"
$ws.Range('E34').Value = "'  +6.99%  "
$ws.Range('B35').Value = "'Kaspa"
$ws.Range('C35').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('D35').Value = "'0.0908"
$ws.Range('E35').Value = "'  +0.25%  "
$ws.Range('B36').Value = "'InternetComputer(DFINITY)"
$ws.Range('C36').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('D36').Value = "'4.32"
$ws.Range('E36').Value = "'  +7.42%  "
$ws.Range('B38').Value = "'WEMIXToken"
$ws.Range('C38').Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range('D38').Value = "'1.83"
$ws.Range('E38').Value = "'  -3.59%  "
$ws.Range('B39').Value = "'LidoDAOToken"
$ws.Range('C39').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('D39').Value = "'2.26"
$ws.Range('E39').Value = "'  -2.09%  "
$ws.Range('D40').Value = "'1.35"
$ws.Range('E40').Value = "'  +0.59%  "
$ws.Range('D41').Value = "'0.0228"
$ws.Range('E41').Value = "'  +5.48%  "
$ws.Range('D42').Value = "'17.89"
$ws.Range('E42').Value = "'  +12.18%  "
$ws.Range('D43').Value = "'1.17"
$ws.Range('E43').Value = "'  +1.81%  "
$ws.Range('B44').Value = "'Cronos"
$ws.Range('C44').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D44').Value = "'0.0973"
$ws.Range('E44').Value = "'  +17.86%  "
$ws.Range('B45').Value = "'Aave"
$ws.Range('C45').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('D45').Value = "'99.53"
$ws.Range('E45').Value = "'  +1.54%  "
$ws.Range('B46').Value = "'HuobiToken"
$ws.Range('C46').Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range('D46').Value = "'2.80"
$ws.Range('E46').Value = "'  +0.46%  "
$ws.Range('B47').Value = "'FTXToken"
$ws.Range('C47').Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range('D47').Value = "'4.15"
$ws.Range('E47').Value = "'  +109.14%  "
$ws.Range('D48').Value = "'1.329.44"
$ws.Range('E48').Value = "'  -0.94%  "
$ws.Range('D49').Value = "'2.94"
$ws.Range('E49').Value = "'  +3.84%  "
$ws.Range('D50').Value = "'2.35"
$ws.Range('E50').Value = "'  +5.46%  "
$ws.Range('E51').Value = "'  +11.82%  "
